$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.43179426924891
$ws.Range("C2").Value = 13.58263742246801
$ws.Range("D2").Value = 4.305216562060667
$ws.Range("F2").Value = 21.62799262871058
$ws.Range("G2").Value = 24.30385587107391
$ws.Range("H2").Value = 13.06525759432081
$ws.Range("L2").Value = 10.97018820084291
$ws.Range("M2").Value = 14.53400870557473
$ws.Range("N2").Value = 17.57753040792734
$ws.Range("O2").Value = 19.30677323522859

$ws.Range("B3").Value = 12.95282786870661
$ws.Range("C3").Value = 13.51707341477304
$ws.Range("D3").Value = 4.240537990337661
$ws.Range("F3").Value = 21.60826923625955
$ws.Range("G3").Value = 24.24783314586591
$ws.Range("H3").Value = 13.101836895945
$ws.Range("L3").Value = 10.98555186920003
$ws.Range("M3").Value = 14.44684330783758
$ws.Range("N3").Value = 17.62419582189658
$ws.Range("O3").Value = 19.34708917045494

$ws.Range("B4").Value = 12.65098184228637
$ws.Range("C4").Value = 13.47697032279141
$ws.Range("D4").Value = 4.199773831348776
$ws.Range("F4").Value = 21.60260451283415
$ws.Range("G4").Value = 24.22262260801078
$ws.Range("H4").Value = 13.12662761563202
$ws.Range("L4").Value = 10.99666218704592
$ws.Range("M4").Value = 14.39519580520319
$ws.Range("N4").Value = 17.65461644056258
$ws.Range("O4").Value = 19.37663375213185

$ws.Range("B5").Value = 12.52621218433989
$ws.Range("C5").Value = 13.46067447466043
$ws.Range("D5").Value = 4.182907914479742
$ws.Range("F5").Value = 21.60191941396169
$ws.Range("G5").Value = 24.21466646926865
$ws.Range("H5").Value = 13.13731566094686
$ws.Range("L5").Value = 11.00161160704275
$ws.Range("M5").Value = 14.37463649810155
$ws.Range("N5").Value = 17.66745846214614
$ws.Range("O5").Value = 19.38987527862455

$ws.Range("B6").Value = 12.50539347067348
$ws.Range("C6").Value = 13.45797159784239
$ws.Range("D6").Value = 4.180092292781135
$ws.Range("F6").Value = 21.60190374385432
$ws.Range("G6").Value = 24.21348546466941
$ws.Range("H6").Value = 13.13912575800396
$ws.Range("L6").Value = 11.00245893894096
$ws.Range("M6").Value = 14.37125255545913
$ws.Range("N6").Value = 17.66961779405719
$ws.Range("O6").Value = 19.39214652761123

$ws.Range("B7").Value = 12.64930604215713
$ws.Range("C7").Value = 13.47675035271
$ws.Range("D7").Value = 4.199547387246223
$ws.Range("F7").Value = 21.6025886983255
$ws.Range("G7").Value = 24.22250591890934
$ws.Range("H7").Value = 13.12676938785424
$ws.Range("L7").Value = 10.99672722829356
$ws.Range("M7").Value = 14.39491653998716
$ws.Range("N7").Value = 17.65478782815362
$ws.Range("O7").Value = 19.37680746952815

$ws.Range("B8").Value = 13.26836302228166
$ws.Range("C8").Value = 13.56000385974692
$ws.Range("D8").Value = 4.283139559332475
$ws.Range("F8").Value = 21.61985620218933
$ws.Range("G8").Value = 24.28263754312263
$ws.Range("H8").Value = 13.07738607718797
$ws.Range("L8").Value = 10.97513767124683
$ws.Range("M8").Value = 14.50357564203986
$ws.Range("N8").Value = 17.59325421828582
$ws.Range("O8").Value = 19.31967837776421

$ws.Range("B9").Value = 14.4135026805117
$ws.Range("C9").Value = 13.72413323200471
$ws.Range("D9").Value = 4.438295802949426
$ws.Range("F9").Value = 21.70468235104631
$ws.Range("G9").Value = 24.47299391399171
$ws.Range("H9").Value = 12.99906112215563
$ws.Range("L9").Value = 10.9460965334259
$ws.Range("M9").Value = 14.73075241185754
$ws.Range("N9").Value = 17.48657703163726
$ws.Range("O9").Value = 19.245766299482

$ws.Range("B10").Value = 15.20435978947502
$ws.Range("C10").Value = 13.84478117558084
$ws.Range("D10").Value = 4.546401939023961
$ws.Range("F10").Value = 21.79773904020828
$ws.Range("G10").Value = 24.65617382385413
$ws.Range("H10").Value = 12.95283067163931
$ws.Range("L10").Value = 10.93284943289858
$ws.Range("M10").Value = 14.90519736417835
$ws.Range("N10").Value = 17.41667882659701
$ws.Range("O10").Value = 19.21482747069713

$ws.Range("B11").Value = 15.551704935392
$ws.Range("C11").Value = 13.89957979596806
$ws.Range("D11").Value = 4.594197286489639
$ws.Range("F11").Value = 21.84664381446996
$ws.Range("G11").Value = 24.7486708966279
$ws.Range("H11").Value = 12.93426125951016
$ws.Range("L11").Value = 10.92857467028239
$ws.Range("M11").Value = 14.98595924269254
$ws.Range("N11").Value = 17.38671055536007
$ws.Range("O11").Value = 19.20584498540103

$ws.Range("B12").Value = 15.68134625190217
$ws.Range("C12").Value = 13.92030977446887
$ws.Range("D12").Value = 4.612088918932745
$ws.Range("F12").Value = 21.86609685052719
$ws.Range("G12").Value = 24.78499031919815
$ws.Range("H12").Value = 12.92758375785051
$ws.Range("L12").Value = 10.92720722151697
$ws.Range("M12").Value = 15.01672353469107
$ws.Range("N12").Value = 17.37562456178423
$ws.Range("O12").Value = 19.20317656311027

$ws.Range("B13").Value = 15.65351129484574
$ws.Range("C13").Value = 13.91584628135872
$ws.Range("D13").Value = 4.608244993357223
$ws.Range("F13").Value = 21.86186595311964
$ws.Range("G13").Value = 24.77711122064018
$ws.Range("H13").Value = 12.92900611191339
$ws.Range("L13").Value = 10.92749055841647
$ws.Range("M13").Value = 15.01009015218619
$ws.Range("N13").Value = 17.37800047110341
$ws.Range("O13").Value = 19.20371864401328

$ws.Range("B14").Value = 15.5624090512358
$ws.Range("C14").Value = 13.90128573697643
$ws.Range("D14").Value = 4.595673446698468
$ws.Range("F14").Value = 21.84822558149473
$ws.Range("G14").Value = 24.7516331640884
$ws.Range("H14").Value = 12.93370479377937
$ws.Range("L14").Value = 10.9284571361839
$ws.Range("M14").Value = 14.98848672072583
$ws.Range("N14").Value = 17.38579324987402
$ws.Range("O14").Value = 19.20561075792392

$ws.Range("B15").Value = 15.5063571812105
$ws.Range("C15").Value = 13.89236397531263
$ws.Range("D15").Value = 4.587945750487664
$ws.Range("F15").Value = 21.83999171546042
$ws.Range("G15").Value = 24.73619467279636
$ws.Range("H15").Value = 12.93662902860726
$ws.Range("L15").Value = 10.92908190367987
$ws.Range("M15").Value = 14.97527701737298
$ws.Range("N15").Value = 17.39060069782323
$ws.Range("O15").Value = 19.20686521507224

$ws.Range("B16").Value = 15.18140078411658
$ws.Range("C16").Value = 13.84119749261062
$ws.Range("D16").Value = 4.543249832177326
$ws.Range("F16").Value = 21.7946743276058
$ws.Range("G16").Value = 24.65031128356888
$ws.Range("H16").Value = 12.95409378741049
$ws.Range("L16").Value = 10.93316399728494
$ws.Range("M16").Value = 14.89994591748896
$ws.Range("N16").Value = 17.41867406846748
$ws.Range("O16").Value = 19.21551704010553

$ws.Range("B17").Value = 14.97879068814776
$ws.Range("C17").Value = 13.80978049394006
$ws.Range("D17").Value = 4.515469935436144
$ws.Range("F17").Value = 21.76854946318313
$ws.Range("G17").Value = 24.59995541912848
$ws.Range("H17").Value = 12.96543851125359
$ws.Range("L17").Value = 10.93611644968125
$ws.Range("M17").Value = 14.85407822145727
$ws.Range("N17").Value = 17.43636410675873
$ws.Range("O17").Value = 19.22212952745387

$ws.Range("B18").Value = 14.86109108467789
$ws.Range("C18").Value = 13.79170282471196
$ws.Range("D18").Value = 4.49936201560028
$ws.Range("F18").Value = 21.7541425688648
$ws.Range("G18").Value = 24.57185667031643
$ws.Range("H18").Value = 12.97219531921685
$ws.Range("L18").Value = 10.93797951063955
$ws.Range("M18").Value = 14.82782992978572
$ws.Range("N18").Value = 17.44671111837852
$ws.Range("O18").Value = 19.22641209048679

$ws.Range("B19").Value = 14.82104358087754
$ws.Range("C19").Value = 13.78558102872807
$ws.Range("D19").Value = 4.493886151264719
$ws.Range("F19").Value = 21.74937135101252
$ws.Range("G19").Value = 24.56249214493127
$ws.Range("H19").Value = 12.97452282567169
$ws.Range("L19").Value = 10.93863864591145
$ws.Range("M19").Value = 14.81896627935859
$ws.Range("N19").Value = 17.45024402867048
$ws.Range("O19").Value = 19.22794436827681

$ws.Range("B20").Value = 15.00048011230276
$ws.Range("C20").Value = 13.8131257216733
$ws.Range("D20").Value = 4.518440635705705
$ws.Range("F20").Value = 21.77126646740973
$ws.Range("G20").Value = 24.6052265783877
$ws.Range("H20").Value = 12.96420687061705
$ws.Range("L20").Value = 10.93578509440008
$ws.Range("M20").Value = 14.8589472383824
$ws.Range("N20").Value = 17.43446315774467
$ws.Range("O20").Value = 19.22137601016429

$ws.Range("B21").Value = 15.58922007256312
$ws.Range("C21").Value = 13.90556316450417
$ws.Range("D21").Value = 4.599371710551385
$ws.Range("F21").Value = 21.85220684198858
$ws.Range("G21").Value = 24.75908181823003
$ws.Range("H21").Value = 12.93231505733618
$ws.Range("L21").Value = 10.92816641306121
$ws.Range("M21").Value = 14.99482741159809
$ws.Range("N21").Value = 17.38349720697262
$ws.Range("O21").Value = 19.20503509877306

$ws.Range("B22").Value = 15.96293541423297
$ws.Range("C22").Value = 13.96584944088234
$ws.Range("D22").Value = 4.651052123792156
$ws.Range("F22").Value = 21.91054303181662
$ws.Range("G22").Value = 24.86715668700552
$ws.Range("H22").Value = 12.91353739612285
$ws.Range("L22").Value = 10.92465168897465
$ws.Range("M22").Value = 15.08468127299906
$ws.Range("N22").Value = 17.35171685117855
$ws.Range("O22").Value = 19.19862822746504

$ws.Range("B23").Value = 15.76451977120753
$ws.Range("C23").Value = 13.93368807255841
$ws.Range("D23").Value = 4.623583003914294
$ws.Range("F23").Value = 21.87891453783375
$ws.Range("G23").Value = 24.8087960174416
$ws.Range("H23").Value = 12.92337026478439
$ws.Range("L23").Value = 10.92639375881871
$ws.Range("M23").Value = 15.03663553759586
$ws.Range("N23").Value = 17.36853894108884
$ws.Range("O23").Value = 19.20165655911723

$ws.Range("B24").Value = 14.99067810422857
$ws.Range("C24").Value = 13.81161339218884
$ws.Range("D24").Value = 4.517098007830849
$ws.Range("F24").Value = 21.77003620064072
$ws.Range("G24").Value = 24.60284083300589
$ws.Range("H24").Value = 12.96476296469065
$ws.Range("L24").Value = 10.93593438408211
$ws.Range("M24").Value = 14.85674557530477
$ws.Range("N24").Value = 17.43532202610431
$ws.Range("O24").Value = 19.22171517725333

$ws.Range("B25").Value = 14.11202461665714
$ws.Range("C25").Value = 13.67968681929122
$ws.Range("D25").Value = 4.39731811893961
$ws.Range("F25").Value = 21.67630700112449
$ws.Range("G25").Value = 24.41381611829822
$ws.Range("H25").Value = 13.01826489542123
$ws.Range("L25").Value = 10.95253087981508
$ws.Range("M25").Value = 14.66789286167335
$ws.Range("N25").Value = 19.26166604865368
